# BOM workbook update: PCB history and stuff, add fab files.
# - Adds three newly-sourced parts (Power NMOS, 120VAC->12VDC converter,
#   test points) as rows 17-19.
# - The diode row's price isn't known yet, so its price cell is blanked
#   out (a stray space) instead of "0", which turns the component total
#   and the assembly total into errors until it's filled in.
# - Clears out the other leftover "0" price placeholders.
# - Leaves the selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New parts discovered while writing up the PCB history / fab files.
$ws.Range("B17").Value = "Power NMOS "
$ws.Range("D17").Value = "LM395T/NOPB-ND"
$ws.Range("G17").Value = "https://www.digikey.com/en/products/detail/texas-instruments/LM395T-NOPB/8902"
$ws.Range("C17").Value = "Digikey"
$ws.Range("E17").Value = 1
$ws.Range("H17").Formula = "=E17*F17"

$ws.Range("B18").Value = "120VAC to 12VDC Converter"
$ws.Range("G18").Value = "https://www.amazon.com/ALITOVE-Converter-Cigarette-Transformer-Refrigerator/dp/B078RZQ9WY/"
$ws.Range("C18").Value = "Amazon"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = 1
$ws.Range("H18").Formula = "=E18*F18"

$ws.Range("D19").Value = "RCU-0C"
$ws.Range("B19").Value = "test points "
$ws.Range("G19").Value = "https://www.digikey.com/en/products/detail/te-connectivity-amp-connectors/RCU-0C/2366048"
$ws.Range("C19").Value = "Digikey"
$ws.Range("E19").Value = 2
$ws.Range("H19").Formula = "=E19*F19"

# 2. Diode (row 7) price isn't known yet -- blank it out instead of "0".
$ws.Range("F7").Value = " "

# 3. Hyperlink the diode's vendor link cell (style picks up the Hyperlink
#    look automatically, same as the other linked cells).
$ws.Hyperlinks.Add($ws.Range("G7"), "https://www.digikey.com/en/products/detail/texas-instruments/LM395T-NOPB/8902")

# 4. Remove the other leftover placeholder "0" price cells.
$priceCellsToClear = @("F5","F6","F8","F9","F10","F11","F12","F13","F14","F15","F16")
foreach ($ref in $priceCellsToClear) {
    $ws.Range($ref).Clear()
}

# 5. Move the active selection like the author left it.
$ws.Range("J12").Select() | Out-Null
